$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6589
$ws.Range("I40").Value = 3481.1428
$ws.Range("K40").Value = 3481.1428
$ws.Range("M40").Value = -3306.1428

$ws.Range("H53").Value = 758.0769
$ws.Range("J53").Value = 1434.5
$ws.Range("L53").Value = 1434.5
$ws.Range("N53").Value = -2708.5

$ws.Range("H62").Value = 4018
$ws.Range("I62").Value = 3740.1428
$ws.Range("J62").Value = 4666.3335
$ws.Range("K62").Value = 3740.1428
$ws.Range("L62").Value = 4666.3335
$ws.Range("M62").Value = -3116.1428
$ws.Range("N62").Value = -5914.3335

$ws.Range("H65").Value = 4018
$ws.Range("I65").Value = 3740.1428
$ws.Range("J65").Value = 4666.3335
$ws.Range("K65").Value = 18700.714
$ws.Range("L65").Value = 23331.6675
$ws.Range("M65").Value = -15580.714
$ws.Range("N65").Value = -29571.6675

$ws.Range("H98").Value = 425.5625
$ws.Range("I98").Value = 425.5625
$ws.Range("K98").Value = 425.5625
$ws.Range("M98").Value = 1072.4375

$ws.Range("H107").Value = 126
$ws.Range("I107").Value = 129.71428
$ws.Range("K107").Value = 129.71428
$ws.Range("M107").Value = 1790.28572

$ws.Range("H122").Value = 425.5625
$ws.Range("I122").Value = 425.5625
$ws.Range("K122").Value = 1276.6875
$ws.Range("M122").Value = 1173.3125

$ws.Range("H132").Value = 669.1967
$ws.Range("I132").Value = 676.678
$ws.Range("K132").Value = 2030.034
$ws.Range("M132").Value = 499.9659999999999

$ws.Range("H138").Value = 4084.6538
$ws.Range("J138").Value = 4725.75
$ws.Range("L138").Value = 14177.25
$ws.Range("N138").Value = -24457.25

$ws.Range("H141").Value = 1871.8918
$ws.Range("I141").Value = 1871.8918
$ws.Range("K141").Value = 5615.6754
$ws.Range("M141").Value = -435.6754000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4906.706
$ws.Range("I61").Value = 2551.875
$ws.Range("J61").Value = 6999.8887
$ws.Range("K61").Value = 2551.875
$ws.Range("L61").Value = 6999.8887
$ws.Range("M61").Value = -2339.875
$ws.Range("N61").Value = -7423.8887

$ws.Range("H74").Value = 4142.6665
$ws.Range("I74").Value = 3864.8333
$ws.Range("J74").Value = 4698.3335
$ws.Range("K74").Value = 3864.8333
$ws.Range("L74").Value = 4698.3335
$ws.Range("M74").Value = -2990.8333
$ws.Range("N74").Value = -6446.3335

$ws.Range("H77").Value = 4142.6665
$ws.Range("I77").Value = 3864.8333
$ws.Range("J77").Value = 4698.3335
$ws.Range("K77").Value = 19324.1665
$ws.Range("L77").Value = 23491.6675
$ws.Range("M77").Value = -14956.1665
$ws.Range("N77").Value = -32227.6675

$ws.Range("H136").Value = 4906.706
$ws.Range("I136").Value = 2551.875
$ws.Range("J136").Value = 6999.8887
$ws.Range("K136").Value = 7655.625
$ws.Range("L136").Value = 20999.6661
$ws.Range("M136").Value = -5105.625
$ws.Range("N136").Value = -26099.6661

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3714.8823
$ws.Range("I86").Value = 3049.75
$ws.Range("K86").Value = 3049.75
$ws.Range("M86").Value = -1926.75

$ws.Range("H89").Value = 3714.8823
$ws.Range("I89").Value = 3049.75
$ws.Range("K89").Value = 15248.75
$ws.Range("M89").Value = -9632.75

$ws.Range("H107").Value = 551.9375
$ws.Range("I107").Value = 534.5
$ws.Range("K107").Value = 534.5
$ws.Range("M107").Value = 1385.5

$ws.Range("H134").Value = 2779.6978
$ws.Range("I134").Value = 1763.6857
$ws.Range("K134").Value = 5291.0571
$ws.Range("M134").Value = -2756.0571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 933.3333
$ws.Range("J23").Value = 933.3333
$ws.Range("L23").Value = 933.3333
$ws.Range("N23").Value = -1413.3333

$ws.Range("H27").Value = 933.3333
$ws.Range("J27").Value = 933.3333
$ws.Range("L27").Value = 933.3333
$ws.Range("N27").Value = -1317.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4433.3213
$ws.Range("J5").Value = 9894.416999999999
$ws.Range("L5").Value = 29683.251
$ws.Range("N5").Value = -29907.251

$ws.Range("H97").Value = 1365
$ws.Range("I97").Value = 933.3333
$ws.Range("J97").Value = 2012.5
$ws.Range("K97").Value = 2799.9999
$ws.Range("L97").Value = 6037.5
$ws.Range("M97").Value = -2303.9999
$ws.Range("N97").Value = -7029.5

$ws.Range("H130").Value = 3256.75
$ws.Range("J130").Value = 7000
$ws.Range("L130").Value = 21000
$ws.Range("N130").Value = -31040

$ws.Range("H135").Value = 4433.3213
$ws.Range("J135").Value = 9894.416999999999
$ws.Range("L135").Value = 89049.753
$ws.Range("N135").Value = -94119.753

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 39577.637
$ws.Range("I70").Value = 67267
$ws.Range("K70").Value = 67267
$ws.Range("M70").Value = -66997

$ws.Range("H73").Value = 39577.637
$ws.Range("I73").Value = 67267
$ws.Range("K73").Value = 67267
$ws.Range("M73").Value = -66331

$ws.Range("H80").Value = 22299568
$ws.Range("I80").Value = 104107.63
$ws.Range("J80").Value = 83337090
$ws.Range("K80").Value = 104107.63
$ws.Range("L80").Value = 83337090
$ws.Range("M80").Value = -103109.63
$ws.Range("N80").Value = -83339086

$ws.Range("H83").Value = 22299568
$ws.Range("I83").Value = 104107.63
$ws.Range("J83").Value = 83337090
$ws.Range("K83").Value = 520538.15
$ws.Range("L83").Value = 416685450
$ws.Range("M83").Value = -515546.15
$ws.Range("N83").Value = -416695434

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 299989.66
$ws.Range("J116").Value = 299989.66
$ws.Range("L116").Value = 299989.66
$ws.Range("N116").Value = -309167.66

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 18800
$ws.Range("J18").Value = 18800
$ws.Range("L18").Value = 18800
$ws.Range("N18").Value = -19146

$ws.Range("H24").Value = 23800
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 23800
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 23800
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -24260

$ws.Range("H74").Value = 13415.667
$ws.Range("J74").Value = 15126
$ws.Range("L74").Value = 15126
$ws.Range("N74").Value = -16998

$ws.Range("H77").Value = 13415.667
$ws.Range("J77").Value = 15126
$ws.Range("L77").Value = 45378
$ws.Range("N77").Value = -54738

$ws.Range("H122").Value = 4212.3125
$ws.Range("I122").Value = 2672
$ws.Range("K122").Value = 8016
$ws.Range("M122").Value = -5566

$ws.Range("H132").Value = 2879.7812
$ws.Range("I132").Value = 1485.1578
$ws.Range("J132").Value = 4918.077
$ws.Range("K132").Value = 4455.4734
$ws.Range("L132").Value = 14754.231
$ws.Range("M132").Value = -1925.4734
$ws.Range("N132").Value = -19814.231
